$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title
Replace-Text "Quantum Entanglement Unveiled" "The Symphony of Chemistry: A Molecular Dance of Transformation"

# Author name "Brian Weiss" -> "Dr. Elijah Brooks"
Replace-Text "Brian Weiss" "Dr. Elijah Brooks"

# Email parts
Replace-Text "brian" "elijah"
Replace-Text "weiss@quantumstudies" "brooks@educonnect"
Replace-Text "com" "org"

# Body paragraph 1
Replace-Text "In the enigmatic realm of quantum physics, where the very fabric of reality appears to challenge our conventional understanding, lies a phenomenon that has captivated the imaginations of scientists and philosophers alike: quantum entanglement." "In the vast orchestra of science, chemistry stands as a harmonious symphony of elements, compounds, and reactions."

Replace-Text " This peculiar behavior, where two particles become linked in such a manner that the state of one instantaneously influences the state of the other, regardless of the distance between them, has baffled researchers for decades." " This intricate dance of particles, governed by the laws of nature, unveils the secrets of matter and its transformations."

Replace-Text " It has ignited a multitude of experiments, sparked heated debates, and opened up new avenues for exploration in the realm of quantum information processing. Delve into the perplexities of quantum entanglement, as we unravel its profound implications for our comprehension of the universe." " As we delve into the realm of chemistry, we embark on a journey through the microscopic world, where atoms and molecules engage in a symphony of change, revealing the fundamental principles that shape our universe."

Replace-Text "Our journey into the realm of quantum entanglement begins with the discovery that certain particles, such as electrons, can be correlated in a way that defies classical intuition." "Like a conductor orchestrating a symphony, chemists manipulate atoms and molecules, creating new substances with unique properties."

Replace-Text " These pairs of entangled particles exhibit an exquisite and puzzling synchronicity." " They unlock the secrets of chemical reactions, harnessing their power to synthesize medicines, fuels, and materials that enhance our lives."

Replace-Text " If the spin state of one particle is measured, the spin state of its entangled partner instantaneously adjusts in accordance, regardless of whether they are separated by a mere nanometer or vast lightyears. This non-local connection between entangled particles violates the principle of locality--a cornerstone of classical physics--and has led to spirited discussions regarding the nature of reality itself." " Chemistry is the language of the natural world, a symphony of particles that echoes through the universe, connecting the smallest particles of matter to the grandest phenomena of life."

Replace-Text "The implications of quantum entanglement extend far beyond the realm of theoretical physics." "Chemistry is not merely a collection of abstract concepts; it is a living tapestry woven into the fabric of our everyday existence."

Replace-Text " This phenomenon has ignited a revolution in the field of quantum information processing." " From the air we breathe to the food we eat, chemistry plays a vital role in sustaining life."

Replace-Text " Entangled particles hold the potential to revolutionize fields such as cryptography, where the security of communication rests upon the intrinsic randomness and nonlocal nature of quantum physics." " The clothes we wear, the medicines that heal us, and the technologies that connect us all bear the imprint of chemistry's transformative power."

Replace-Text " Furthermore, entanglement-based quantum computing holds the promise of unprecedented computational power, dwarfing the capabilities of today's most advanced supercomputers. As we venture further into this enigmatic domain, the potential applications of quantum entanglement appear boundless, beckoning us to push the very limits of our scientific and technological frontiers." " It is a science that bridges the gap between the natural and the artificial, shaping our world in countless ways."

# Summary paragraph
Replace-Text "Quantum entanglement stands as a testament to the interconnectedness of the universe and challenges our conventional notions of space, time, and causality." "Through the lens of chemistry, we gain a profound understanding of the world around us."

Replace-Text " Through a series of intricate experiments, scientists have demonstrated that entangled particles exhibit an uncanny link, where the state of one instantly influences the state of the other, regardless of the distance between them." " Chemistry is a symphony of particles, a dance of atoms and molecules that reveals the secrets of matter and its transformations."

Replace-Text " This non-local connection has ignited fervent debates and spurred scientific exploration into the nature of reality." " It is a science that transcends the boundaries of abstraction, impacting our daily lives in countless ways."

Replace-Text " Furthermore, the harnessing of quantum entanglement holds immense promise for advancements in fields such as quantum cryptography, quantum computing, and teleportation. As we delve deeper into the enigmatic realm of quantum entanglement, we unlock new avenues for understanding the universe and pave the way for groundbreaking technologies that will shape the future of human civilization." " Chemistry is the language of nature, the conductor of life's grand symphony."

Write-Output "Done"
